# Update the "想去人数" (want-to-go count) figures for two events that
# each appear in both the "展览" sheet and the "全部类型" sheet.
#   F2: 367 -> 368
#   F10: 429 -> 430

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 368
    $ws.Range("F10").Value = 430
}
